$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings) in row 1
$ws.Range("B1").Value = "B_LF"
$ws.Range("C1").Value = "C_B"
$ws.Range("D1").Value = "C_FFR"
$ws.Range("E1").Value = "LF_C"

# Update numeric values in row 2
$ws.Range("B2").Value = 0.01139498042573416
$ws.Range("C2").Value = 0.9400827977904548
$ws.Range("D2").Value = 0.001956093852152114
$ws.Range("E2").Value = 11.90697295111552

# Update numeric values in row 3
$ws.Range("B3").Value = 0.04638754707348047
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.02023979304397994
$ws.Range("E3").Value = 0.006868588110146678
